$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.117.86"
$ws.Range("E2").Value = "  -0.85%  "

$ws.Range("D3").Value = "1.823.46"
$ws.Range("E3").Value = "  -1.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.22"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4624"
$ws.Range("E7").Value = "  -2.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3635"
$ws.Range("E8").Value = "  -1.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07300"
$ws.Range("E9").Value = "  -2.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8697"
$ws.Range("E10").Value = "  -1.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.11"
$ws.Range("E11").Value = "  -1.95%  "

$ws.Range("D12").Value = "1.864.21"
$ws.Range("E12").Value = "  +1.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07581"
$ws.Range("E13").Value = "  +3.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.342"
$ws.Range("E14").Value = "  -2.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.41"
$ws.Range("E15").Value = "  -1.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.495"
$ws.Range("E16").Value = "  -1.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008637"
$ws.Range("E18").Value = "  -2.41%  "

$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("D20").Value = "27.426.05"
$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.47"
$ws.Range("E21").Value = "  -2.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.206"
$ws.Range("E22").Value = "  -2.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.56"
$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("D24").Value = "2.095.64"
$ws.Range("E24").Value = "  +1.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.94"
$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.871"
$ws.Range("E26").Value = "  -1.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.23"
$ws.Range("E27").Value = "  -2.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.078"
$ws.Range("E28").Value = "  -5.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.095"
$ws.Range("E29").Value = "  -3.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.14"
$ws.Range("E30").Value = "  -1.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08902"
$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("E32").Value = "  +0.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7316"
$ws.Range("E33").Value = "  -4.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.455"
$ws.Range("E34").Value = "  -2.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.137"
$ws.Range("E35").Value = "  -3.62%  "

$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.458"
$ws.Range("E37").Value = "  +2.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.072"
$ws.Range("E38").Value = "  -3.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05246"
$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01915"
$ws.Range("E40").Value = "  -2.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.925"
$ws.Range("E41").Value = "  -2.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.130"
$ws.Range("E42").Value = "  -3.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5189"
$ws.Range("E43").Value = "  -3.36%  "

$ws.Range("E44").Value = "  -2.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.260"
$ws.Range("E45").Value = "  -3.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4844"
$ws.Range("E46").Value = "  -2.75%  "

$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.22"
$ws.Range("E48").Value = "  -3.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.48"
$ws.Range("E49").Value = "  -0.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.633"
$ws.Range("E50").Value = "  -2.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06264"
$ws.Range("E51").Value = "  -0.93%  "

